$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("L3").Value2 = 0
$ws.Range("N3").Value2 = $null
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 0
$ws.Range("K16").Value2 = 0
$ws.Range("M16").Value2 = $null
$ws.Range("H64").Value2 = 3400
$ws.Range("J64").Value2 = 3400
$ws.Range("L64").Value2 = 3400
$ws.Range("N64").Value2 = -3896
$ws.Range("H67").Value2 = 3400
$ws.Range("J67").Value2 = 3400
$ws.Range("L67").Value2 = 3400
$ws.Range("N67").Value2 = -5116
$ws.Range("H102").Value2 = 0
$ws.Range("J102").Value2 = 0
$ws.Range("L102").Value2 = 0
$ws.Range("N102").Value2 = $null
$ws.Range("H111").Value2 = 2164.8333
$ws.Range("I111").Value2 = 3937.6
$ws.Range("K111").Value2 = 11812.8
$ws.Range("M111").Value2 = -8745.799999999999
$ws.Range("H129").Value2 = 755.89655
$ws.Range("J129").Value2 = 799.2075
$ws.Range("L129").Value2 = 2397.6225
$ws.Range("N129").Value2 = -12397.6225
$ws.Range("H132").Value2 = 5507.8335
$ws.Range("I132").Value2 = 5452.5625
$ws.Range("J132").Value2 = 5950
$ws.Range("K132").Value2 = 16357.6875
$ws.Range("L132").Value2 = 17850
$ws.Range("M132").Value2 = -13827.6875
$ws.Range("N132").Value2 = -22910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1176.9412
$ws.Range("I2").Value2 = 1108.3077
$ws.Range("K2").Value2 = 1108.3077
$ws.Range("M2").Value2 = -995.3077000000001
$ws.Range("H32").Value2 = 6481.3145
$ws.Range("I32").Value2 = 4893.923
$ws.Range("J32").Value2 = 17737.363
$ws.Range("K32").Value2 = 4893.923
$ws.Range("L32").Value2 = 17737.363
$ws.Range("M32").Value2 = -4606.923
$ws.Range("N32").Value2 = -18311.363
$ws.Range("H45").Value2 = 2350
$ws.Range("I45").Value2 = 2217.1765
$ws.Range("J45").Value2 = 2575.8
$ws.Range("K45").Value2 = 2217.1765
$ws.Range("L45").Value2 = 2575.8
$ws.Range("M45").Value2 = -1840.1765
$ws.Range("N45").Value2 = -3329.8
$ws.Range("H61").Value2 = 1647.711
$ws.Range("I61").Value2 = 1378.5428
$ws.Range("J61").Value2 = 2589.8
$ws.Range("K61").Value2 = 1378.5428
$ws.Range("L61").Value2 = 2589.8
$ws.Range("M61").Value2 = -1166.5428
$ws.Range("N61").Value2 = -3013.8
$ws.Range("H110").Value2 = 1218.625
$ws.Range("I110").Value2 = 1170.5714
$ws.Range("J110").Value2 = 1555
$ws.Range("K110").Value2 = 1170.5714
$ws.Range("L110").Value2 = 1555
$ws.Range("M110").Value2 = 874.4286
$ws.Range("N110").Value2 = -5645
$ws.Range("H116").Value2 = 1176.9412
$ws.Range("I116").Value2 = 1108.3077
$ws.Range("K116").Value2 = 1108.3077
$ws.Range("M116").Value2 = 1185.6923
$ws.Range("H122").Value2 = 2397.7778
$ws.Range("I122").Value2 = 2185
$ws.Range("K122").Value2 = 6555
$ws.Range("M122").Value2 = -4105
$ws.Range("H132").Value2 = 12684.596
$ws.Range("I132").Value2 = 1907.2572
$ws.Range("J132").Value2 = 44118.5
$ws.Range("K132").Value2 = 5721.7716
$ws.Range("L132").Value2 = 132355.5
$ws.Range("M132").Value2 = -3191.7716
$ws.Range("N132").Value2 = -137415.5
$ws.Range("H136").Value2 = 1647.711
$ws.Range("I136").Value2 = 1378.5428
$ws.Range("J136").Value2 = 2589.8
$ws.Range("K136").Value2 = 4135.6284
$ws.Range("L136").Value2 = 7769.400000000001
$ws.Range("M136").Value2 = -1585.6284
$ws.Range("N136").Value2 = -12869.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1176.9412
$ws.Range("I3").Value2 = 1108.3077
$ws.Range("K3").Value2 = 1108.3077
$ws.Range("M3").Value2 = -994.3077000000001
$ws.Range("H94").Value2 = 765.7273
$ws.Range("I94").Value2 = 652.25
$ws.Range("J94").Value2 = 1068.3334
$ws.Range("K94").Value2 = 652.25
$ws.Range("L94").Value2 = 1068.3334
$ws.Range("M94").Value2 = -201.25
$ws.Range("N94").Value2 = -1970.3334
$ws.Range("H105").Value2 = 1138267.5
$ws.Range("I105").Value2 = 1561.5264
$ws.Range("J105").Value2 = 2002164
$ws.Range("K105").Value2 = 1561.5264
$ws.Range("L105").Value2 = 2002164
$ws.Range("M105").Value2 = 185.4736
$ws.Range("N105").Value2 = -2005658
$ws.Range("H107").Value2 = 1497.963
$ws.Range("I107").Value2 = 1461.1364
$ws.Range("J107").Value2 = 1660
$ws.Range("K107").Value2 = 1461.1364
$ws.Range("L107").Value2 = 1660
$ws.Range("M107").Value2 = 458.8635999999999
$ws.Range("N107").Value2 = -5500
$ws.Range("H134").Value2 = 4180.2334
$ws.Range("I134").Value2 = 4180.2334
$ws.Range("K134").Value2 = 12540.7002
$ws.Range("M134").Value2 = -10005.7002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1306.1
$ws.Range("I16").Value2 = 1231
$ws.Range("J16").Value2 = 1338.2858
$ws.Range("K16").Value2 = 1231
$ws.Range("L16").Value2 = 1338.2858
$ws.Range("M16").Value2 = -944
$ws.Range("N16").Value2 = -1912.2858
$ws.Range("H31").Value2 = 3527.39
$ws.Range("I31").Value2 = 1758.8518
$ws.Range("K31").Value2 = 1758.8518
$ws.Range("M31").Value2 = -1463.8518
$ws.Range("H34").Value2 = 3527.39
$ws.Range("I34").Value2 = 1758.8518
$ws.Range("K34").Value2 = 1758.8518
$ws.Range("M34").Value2 = -1556.8518
$ws.Range("I99").Value2 = 2936.5334
$ws.Range("J99").Value2 = 6400
$ws.Range("K99").Value2 = 2936.5334
$ws.Range("L99").Value2 = 6400
$ws.Range("M99").Value2 = -1438.5334
$ws.Range("N99").Value2 = -9396
$ws.Range("H107").Value2 = 993.80646
$ws.Range("J107").Value2 = 1698.6428
$ws.Range("L107").Value2 = 1698.6428
$ws.Range("N107").Value2 = -5538.6428
$ws.Range("H113").Value2 = 1306.1
$ws.Range("I113").Value2 = 1231
$ws.Range("J113").Value2 = 1338.2858
$ws.Range("K113").Value2 = 1231
$ws.Range("L113").Value2 = 1338.2858
$ws.Range("M113").Value2 = 939
$ws.Range("N113").Value2 = -5678.2858
$ws.Range("I126").Value2 = 2936.5334
$ws.Range("J126").Value2 = 6400
$ws.Range("K126").Value2 = 8809.600199999999
$ws.Range("L126").Value2 = 19200
$ws.Range("M126").Value2 = -6339.600199999999
$ws.Range("N126").Value2 = -24140
$ws.Range("H132").Value2 = 4234.25
$ws.Range("J132").Value2 = 6102
$ws.Range("L132").Value2 = 18306
$ws.Range("N132").Value2 = -23366
$ws.Range("H141").Value2 = 27728.36
$ws.Range("J141").Value2 = 27728.36
$ws.Range("L141").Value2 = 27728.36
$ws.Range("N141").Value2 = -38088.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 753.6070999999999
$ws.Range("I122").Value2 = 548
$ws.Range("J122").Value2 = 798.3043
$ws.Range("K122").Value2 = 4932
$ws.Range("L122").Value2 = 7184.7387
$ws.Range("M122").Value2 = -2482
$ws.Range("N122").Value2 = -12084.7387
$ws.Range("H131").Value2 = 735.28
$ws.Range("J131").Value2 = 748.4639
$ws.Range("L131").Value2 = 2245.3917
$ws.Range("N131").Value2 = -12325.3917
$ws.Range("H137").Value2 = 20840190
$ws.Range("I137").Value2 = 2140
$ws.Range("J137").Value2 = 27786206
$ws.Range("K137").Value2 = 6420
$ws.Range("L137").Value2 = 83358618
$ws.Range("M137").Value2 = -1320
$ws.Range("N137").Value2 = -83368818

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 3217
$ws.Range("I102").Value2 = 2030.75
$ws.Range("K102").Value2 = 2030.75
$ws.Range("M102").Value2 = -408.75
$ws.Range("H126").Value2 = 2938.6604
$ws.Range("I126").Value2 = 2076.5
$ws.Range("J126").Value2 = 3652.1724
$ws.Range("K126").Value2 = 6229.5
$ws.Range("L126").Value2 = 10956.5172
$ws.Range("M126").Value2 = -3759.5
$ws.Range("N126").Value2 = -15896.5172
$ws.Range("H132").Value2 = 16517.764
$ws.Range("I132").Value2 = 3461.7812
$ws.Range("K132").Value2 = 10385.3436
$ws.Range("M132").Value2 = -7855.3436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 3841.7896
$ws.Range("I40").Value2 = 3886.2666
$ws.Range("J40").Value2 = 3675
$ws.Range("K40").Value2 = 3886.2666
$ws.Range("L40").Value2 = 3675
$ws.Range("M40").Value2 = -3750.2666
$ws.Range("N40").Value2 = -3947
$ws.Range("H82").Value2 = 1078.6
$ws.Range("I82").Value2 = 1093.5454
$ws.Range("K82").Value2 = 1093.5454
$ws.Range("M82").Value2 = -732.5454
$ws.Range("H85").Value2 = 1078.6
$ws.Range("I85").Value2 = 1093.5454
$ws.Range("K85").Value2 = 1093.5454
$ws.Range("M85").Value2 = 154.4546
$ws.Range("H100").Value2 = 2128.5386
$ws.Range("I100").Value2 = 797.6667
$ws.Range("K100").Value2 = 797.6667
$ws.Range("M100").Value2 = -256.6667
$ws.Range("H122").Value2 = 2182713.8
$ws.Range("I122").Value2 = 2454590.5
$ws.Range("K122").Value2 = 7363771.5
$ws.Range("M122").Value2 = -7361321.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 1749.5758
$ws.Range("I126").Value2 = 1365.5714
$ws.Range("K126").Value2 = 4096.7142
$ws.Range("M126").Value2 = -1626.7142
$ws.Range("H132").Value2 = 827.8889
$ws.Range("I132").Value2 = 721.3461
$ws.Range("K132").Value2 = 2164.0383
$ws.Range("M132").Value2 = 365.9616999999998
$ws.Range("H136").Value2 = 34412016
$ws.Range("I136").Value2 = 39703730
$ws.Range("J136").Value2 = 15875
$ws.Range("K136").Value2 = 119111190
$ws.Range("L136").Value2 = 47625
$ws.Range("M136").Value2 = -119108640
$ws.Range("N136").Value2 = -52725
